# Fruta / hortaliza, semanal
# Insert two new weekly report rows (rows 57-58) for Chirimoya / Vega Modelo
# de Temuco, pushing the previous rows 57-67 down to 59-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 57:67 down by two rows, opening up a gap at 57:58.
$ws.Rows("57:58").Insert()

# --- New row 57 -----------------------------------------------------------
$ws.Cells.Item(57, 1).Value = 10
$ws.Cells.Item(57, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(57, 3).Value = "La Araucanía"
$ws.Cells.Item(57, 4).Value = 44476
$ws.Cells.Item(57, 5).Value = 9
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100107
$ws.Cells.Item(57, 8).Value = "Otros"
$ws.Cells.Item(57, 9).Value = 100107002
$ws.Cells.Item(57, 10).Value = "Chirimoya"
$ws.Cells.Item(57, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(57, 12).Value = "Especial"
$ws.Cells.Item(57, 13).Value = 50
$ws.Cells.Item(57, 14).Value = 3500
$ws.Cells.Item(57, 15).Value = 3500
$ws.Cells.Item(57, 16).Value = 3500
$ws.Cells.Item(57, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(57, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(57, 19).Value = 3500
$ws.Cells.Item(57, 20).Value = 1

# --- New row 58 -----------------------------------------------------------
$ws.Cells.Item(58, 1).Value = 10
$ws.Cells.Item(58, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(58, 3).Value = "La Araucanía"
$ws.Cells.Item(58, 4).Value = 44476
$ws.Cells.Item(58, 5).Value = 9
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100107
$ws.Cells.Item(58, 8).Value = "Otros"
$ws.Cells.Item(58, 9).Value = 100107002
$ws.Cells.Item(58, 10).Value = "Chirimoya"
$ws.Cells.Item(58, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 150
$ws.Cells.Item(58, 14).Value = 3000
$ws.Cells.Item(58, 15).Value = 3000
$ws.Cells.Item(58, 16).Value = 3000
$ws.Cells.Item(58, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(58, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(58, 19).Value = 3000
$ws.Cells.Item(58, 20).Value = 1
